$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Update the Date value (B8, next to "Date" in A8)
$ws.Range("B8").Value = "2023-02-01T09:05:11-06:00"

# Add the required "Experimental" boolean value (true) as text in B7,
# keeping the cell's existing style/border (s="2") and string type (t="s")
# rather than Excel's native boolean type (t="b").
$ws.Range("B7").Formula = "=""true"""
$ws.Range("B7").Copy()
$ws.Range("B7").PasteSpecial(-4163)
